$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 45236
$ws.Range("B8").Value = 82
$ws.Range("C8").Value = "Added postImage screen, setup backend for posts, tabBar. Created post model and started feed Screen work. "
$ws.Range("D8").Value = 1

$ws.Range("C12").Select()
